$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update gamer names across columns A, D, G, J (rows 2-17)
# Written in the same order as the original authoring to preserve shared-string ordering
$ws.Range("J5").Value2 = "HolyGamer"
$ws.Range("J6").Value2 = "TallGamer"
$ws.Range("J4").Value2 = "SamuraiGamer"
$ws.Range("J7").Value2 = "Not-A-Gamer"
$ws.Range("G5").Value2 = "HighSocietyGamer"
$ws.Range("J2").Value2 = "SocietyGamer"
$ws.Range("G2").Value2 = "HackerGamer"
$ws.Range("G3").Value2 = "FrenchGamer"
$ws.Range("G4").Value2 = "DelayedGamer"
$ws.Range("G9").Value2 = "CursedGamer"
$ws.Range("G8").Value2 = "StillAGamer"
$ws.Range("G7").Value2 = "SweatyGamer"
$ws.Range("G6").Value2 = "OrangeGamer"
$ws.Range("D12").Value2 = "ConfusedGamer"
$ws.Range("D11").Value2 = "CowboyGamer"
$ws.Range("D9").Value2 = "AnurognathusGamer"
$ws.Range("D4").Value2 = "SwirlyGamer"
$ws.Range("D2").Value2 = "SexyGamer"
$ws.Range("D3").Value2 = "SushiGamer"
$ws.Range("D7").Value2 = "RevolutionaryGamer"
$ws.Range("D8").Value2 = "Colonioal Gamer"
$ws.Range("D6").Value2 = "MysteriousGamer"
$ws.Range("D10").Value2 = "Dragonborn Gamer"
$ws.Range("D13").Value2 = "FrogGamer"
$ws.Range("D5").Value2 = "PipeGamer"
$ws.Range("A15").Value2 = "VillainGamer"
$ws.Range("A11").Value2 = "HillbillyGamer"
$ws.Range("A7").Value2 = "MathGamer"
$ws.Range("A16").Value2 = "NoseGodGamer"
$ws.Range("A5").Value2 = "YoungGamer"
$ws.Range("A12").Value2 = "PunishedGamer"
$ws.Range("A9").Value2 = "TrueGamer"
$ws.Range("A2").Value2 = "SociallyDistantGamer"
$ws.Range("A4").Value2 = "NordicGamer"
$ws.Range("A13").Value2 = "OneEyedGamer"
$ws.Range("A3").Value2 = "AbstractGamer"
$ws.Range("A17").Value2 = "SpecialGamer"
$ws.Range("A8").Value2 = "LuckyGamer"
$ws.Range("A10").Value2 = "ConstructionGamer"
$ws.Range("A6").Value2 = "CuteSchoolGamer"
$ws.Range("A14").Value2 = "AverageGamer"
$ws.Range("J3").Value2 = "GoldenGamer"

# Swap the E12/F12 and E13/F13 values
$ws.Range("E12").Value2 = -25
$ws.Range("F12").Value2 = -75
$ws.Range("E13").Value2 = -70
$ws.Range("F13").Value2 = -210

# Update the active selection
$ws.Range("F24").Select()
